$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update conversion text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$cellA1 = $ws1.Range("A1")
$text = $cellA1.Value2
$text = $text -replace [regex]::Escape("1000 Bs = 1.76 = 6397.3 pesos"), "1000 Bs = 1.86 = 6733.56 pesos"
$text = $text -replace [regex]::Escape("6397.3 pesos = 1.75 = 922.85 Bs"), "6733.56 pesos = 1.85 = 887.03 Bs"
$cellA1.Value = $text

# --- Sheet "tasas": update N10, O10, N12, O12 ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 538.5
$ws2.Range("O10").Value = 3626.02
$ws2.Range("N12").Value = 3643.76
$ws2.Range("O12").Value = 480
